$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range('D2').Value = '64.792.26'
$ws.Range('E2').Value = '  -0.29%  '

$ws.Range('D3').Value = '3.442.76'
$ws.Range('E3').Value = '  -0.78%  '

$ws.Range('E4').Value = '  -0.03%  '

Set-TextValue $ws.Range('D5') '573.18'
$ws.Range('E5').Value = '  -0.64%  '

Set-TextValue $ws.Range('D6') '159.32'
$ws.Range('E6').Value = '  -0.44%  '

$ws.Range('E7').Value = '  -0.03%  '

$ws.Range('D8').Value = '3.444.79'
$ws.Range('E8').Value = '  -0.80%  '

Set-TextValue $ws.Range('D9') '0.574'
$ws.Range('E9').Value = '  -5.34%  '

Set-TextValue $ws.Range('D10') '7.20'
$ws.Range('E10').Value = '  -0.75%  '

$ws.Range('E11').Value = '  -2.15%  '

$ws.Range('E12').Value = '  -0.17%  '

$ws.Range('D13').Value = '4.036.81'
$ws.Range('E13').Value = '  -1.00%  '

$ws.Range('E14').Value = '  -0.40%  '

Set-TextValue $ws.Range('D15') '27.75'
$ws.Range('E15').Value = '  -2.50%  '

Set-TextValue $ws.Range('D16') '0.0000176'
$ws.Range('E16').Value = '  -8.14%  '

$ws.Range('D17').Value = '64.823.79'
$ws.Range('E17').Value = '  -0.37%  '

$ws.Range('D18').Value = '3.443.95'
$ws.Range('E18').Value = '  -1.26%  '

Set-TextValue $ws.Range('D19') '6.27'
$ws.Range('E19').Value = '  -2.20%  '

Set-TextValue $ws.Range('D20') '13.87'
$ws.Range('E20').Value = '  -2.78%  '

Set-TextValue $ws.Range('D21') '379.15'
$ws.Range('E21').Value = '  -0.63%  '

Set-TextValue $ws.Range('D22') '7.96'
$ws.Range('E22').Value = '  -2.57%  '

Set-TextValue $ws.Range('D23') '0.542'
$ws.Range('E23').Value = '  -1.17%  '

$ws.Range('E24').Value = '  +0.29%  '

Set-TextValue $ws.Range('D25') '71.98'
$ws.Range('E25').Value = '  -1.17%  '

$ws.Range('E26').Value = '  +0.53%  '

$ws.Range('E27').Value = '  -1.30%  '

Set-TextValue $ws.Range('D28') '0.178'
$ws.Range('E28').Value = '  -0.79%  '

$ws.Range('E29').Value = '  +0.04%  '

Set-TextValue $ws.Range('D30') '1.48'
$ws.Range('E30').Value = '  +0.42%  '

Set-TextValue $ws.Range('D31') '6.10'
$ws.Range('E31').Value = '  -1.35%  '

Set-TextValue $ws.Range('D32') '2.01'
$ws.Range('E32').Value = '  -2.21%  '

Set-TextValue $ws.Range('D33') '23.18'
$ws.Range('E33').Value = '  -1.42%  '

Set-TextValue $ws.Range('D34') '7.08'
$ws.Range('E34').Value = '  -1.99%  '

$ws.Range('E35').Value = '  -0.32%  '

Set-TextValue $ws.Range('D36') '161.11'
$ws.Range('E36').Value = '  -0.09%  '

Set-TextValue $ws.Range('D37') '1.89'
$ws.Range('E37').Value = '  -0.64%  '

$ws.Range('D38').Value = '2.903.07'
$ws.Range('E38').Value = '  -3.35%  '

Set-TextValue $ws.Range('D39') '0.0745'
$ws.Range('E39').Value = '  -3.41%  '

$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range('D40') '6.67'
$ws.Range('E40').Value = '  +3.02%  '

$ws.Range('B41').Value = 'EnergySwap'
$ws.Range('C41').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range('D41') '26.24'
$ws.Range('E41').Value = '  -2.45%  '

Set-TextValue $ws.Range('D42') '4.54'
$ws.Range('E42').Value = '  +0.17%  '

Set-TextValue $ws.Range('D43') '42.92'
$ws.Range('E43').Value = '  +0.95%  '

$ws.Range('B44').Value = 'Mantle'
$ws.Range('C44').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue $ws.Range('D44') '0.779'
$ws.Range('E44').Value = '  +0.16%  '

$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range('D45') '0.0312'
$ws.Range('E45').Value = '  -2.73%  '

Set-TextValue $ws.Range('D46') '26.01'
$ws.Range('E46').Value = '  +2.04%  '

$ws.Range('B47').Value = 'dogwifhat'
$ws.Range('C47').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue $ws.Range('D47') '2.27'
$ws.Range('E47').Value = '  +4.85%  '

$ws.Range('B48').Value = 'ONDO'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
Set-TextValue $ws.Range('D48') '1.08'
$ws.Range('E48').Value = '  -2.03%  '

Set-TextValue $ws.Range('D49') '316.95'
$ws.Range('E49').Value = '  -0.37%  '

Set-TextValue $ws.Range('D50') '6.48'
$ws.Range('E50').Value = '  -3.45%  '

$ws.Range('B51').Value = 'SuiNetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
Set-TextValue $ws.Range('D51') '0.843'
$ws.Range('E51').Value = '  -2.05%  '
